$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4: Identificador (D) and Tipo (E) columns ---
# Row 2 (Juan)
$ws.Range("D2").Value = "juanSL"
$ws.Range("E2").Value = "Entity"

# Row 3 (Luis)
$ws.Range("D3").Value = "usuarioLuis"
$ws.Range("E3").Value = "Person"

# Row 4 (Ana)
$ws.Range("D4").Value = "sensorAna"
$ws.Range("E4").Value = "Sensor"

# --- Remove stray formatted-but-empty cell G4 (no longer present in new layout) ---
$ws.Range("G4").Clear()

# --- Add new data rows 5-7 ---
# Row 5: Juan / usuario
$ws.Range("A5").Value = "Juan"
$ws.Range("B5").Value = "1.0,0.2"
$ws.Range("C5").Value = "juan@uniovi.es"
$ws.Range("D5").Value = "usuarioJuan"
$ws.Range("E5").Value = "Person"

# Row 6: RACE entity
$ws.Range("A6").Value = "RACE"
$ws.Range("B6").Value = "1.123,-2.123"
$ws.Range("C6").Value = "avisos@race.es"
$ws.Range("D6").Value = "usuarioRace"
$ws.Range("E6").Value = "Entity"

# Row 7: Sensor, uses a distinct (underlined) font for the whole row
$ws.Range("A7").Value = "SensorTemperatura-A6-PK27"
$ws.Range("B7").Value = "23.231,123.2"
$ws.Range("C7").Value = "tecnico@copinsa.es"
$ws.Range("D7").Value = "usuarioA6-PK27"
$ws.Range("E7").Value = "Sensor"
$ws.Range("A7:E7").Font.Underline = $true

# --- Column widths: widen Nombre (A) and Identificador (D) columns ---
$ws.Columns.Item(1).ColumnWidth = 29.25
$ws.Columns.Item(4).ColumnWidth = 23.6

# --- Selection moves to the newly added last row ---
$ws.Range("A7:XFD7").Select() | Out-Null
